$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text (string) representation
# instead of being auto-coerced into numbers by Excel (e.g. "1.000" -> 1).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.921.35'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '1.767.92'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '328.58'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.4551'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '42.01'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").Value = '0.07383'
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '20.72'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '7.185'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '1.768.11'
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("D17").Value = '92.54'
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").Value = '0.00001059'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '0.06440'
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '16.96'
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").Value = '5.760'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '27.944.20'
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("D25").Value = '2.100'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '158.47'
$ws.Range("E26").Value = '  -3.58%  '
$ws.Range("D27").Value = '20.17'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '1.973.51'
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("D29").Value = '2.140'
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("D30").Value = '123.87'
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").Value = '1.081'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '0.09183'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").Value = '5.609'
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").Value = '11.83'
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").Value = '0.02278'
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").Value = '0.06114'
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("D38").Value = '0.2091'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").Value = '4.948'
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("D40").Value = '0.6255'
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").Value = '7.794'
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").Value = '13.25'
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").Value = '3.745'
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").Value = '0.5851'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D47").Value = '122.47'
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("D48").Value = '1.929'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").Value = '1.131'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("E51").Value = '  +1.86%  '
